$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "variance" column (F) from 0 to 2 for the specified rows
$rows = @(7, 8, 15, 16, 22, 23, 30, 31, 39, 40, 41, 42, 43, 44, 45)
foreach ($r in $rows) {
    $ws.Range("F$r").Value = 2
}

# Update the sheet view: clear the scrolled topLeftCell and change the active selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I38").Select()
